# [Kadastro App] Yeni kayit eklendi: 2916
#
# Appends one new record row to both the master "Kayitlar" sheet and the
# matching birim sheet ("Erdemli"), mirroring the layout of the existing
# rows (Kayit No, Tarih, Birim, Parsel Sayisi, Is, Personeller) - all
# columns are stored as text, same as every other row in these sheets.

$wb = $excel.ActiveWorkbook

$newRow = @(
    "2916",
    "2025-09-08",
    "Erdemli",
    "2",
    "HAZIR_BEYANNAME",
    "AYHAN KARADAYI (K.Teknisyeni), EMİNE ALANLI KIRCILI (K.Mühendisi)"
)

$sheetNames = @("Kayitlar", "Erdemli")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # First completely empty row right after the current data block.
    $lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
    $targetRow = $lastRow + 1

    for ($col = 1; $col -le $newRow.Length; $col++) {
        $cell = $ws.Cells.Item($targetRow, $col)
        # Force text storage (columns like Kayit No / Parsel Sayisi look
        # numeric but every existing cell in these columns is text) while
        # keeping the default "Normal" style used by the rest of the sheet.
        $cell.Value = "'" + $newRow[$col - 1]
        $cell.Style = "Normal"
    }
}
